# Final commit for Recipe scraping
# - Re-apply formatting on Sheet2 data rows so the engine collapses the
#   (formerly duplicated) wrap-text style back onto a single style record.
# - Append two new worksheets ("Sheet3" / "Sheet4") holding the Allergies
#   and Nut Allergies lists, each with a bold header cell, and make the
#   last one ("Sheet4") the active/selected sheet - matching the source
#   workbook's new state.

$wb = $excel.ActiveWorkbook

# --- Sheet2: normalize styling so the duplicate wrap-text style is no
#     longer referenced by any cell (A2:A16 keep wrapping, just re-applied)
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A2:A16").WrapText = $true

# --- Sheet3: "Allergies" list
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "Sheet3"
$ws3.Range("A1").Value = "Allergies"
$ws3.Range("A1").Font.Bold = $true
$ws3.Range("A2").Value = "milk"
$ws3.Range("A3").Value = "soy"
$ws3.Range("A4").Value = "egg"
$ws3.Range("A5").Value = "sesame"
$ws3.Range("A6").Value = "shellfish"
$ws3.Range("A7").Value = "seafood"
[void]$ws3.Range("A1:A7").Select()

# --- Sheet4: "Nut Allergies" list
$ws4 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "Sheet4"
$ws4.Range("A1").Value = "Nut Allergies"
$ws4.Range("A1").Font.Bold = $true
$ws4.Range("A2").Value = "peanuts"
$ws4.Range("A3").Value = "walnuts"
$ws4.Range("A4").Value = "almond"
$ws4.Range("A5").Value = "hazelnut"
$ws4.Range("A6").Value = "cashew"
$ws4.Range("A7").Value = "pecan"
$ws4.Range("A8").Value = "pistachio"
[void]$ws4.Range("A1:A8").Select()

# Sheet4 is the last sheet added, so it is now the active/selected tab -
# mirrors the workbook's bookViews/activeTab and sheetView/tabSelected change.
$ws4.Activate()
